# Scheduled runner update: refresh cached market-board profit figures
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns H:N)
# across several leve tables, per the latest price-check pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1508.7778
$ws.Range("I46").Value = 1733.3334
$ws.Range("J46").Value = 1396.5
$ws.Range("K46").Value = 5200.0002
$ws.Range("L46").Value = 4189.5
$ws.Range("M46").Value = -5081.0002
$ws.Range("N46").Value = -4427.5
$ws.Range("H60").Value = 1508.7778
$ws.Range("I60").Value = 1733.3334
$ws.Range("J60").Value = 1396.5
$ws.Range("K60").Value = 5200.0002
$ws.Range("L60").Value = 4189.5
$ws.Range("M60").Value = -4716.0002
$ws.Range("N60").Value = -5157.5
$ws.Range("H69").Value = 3497.4
$ws.Range("I69").Value = 2500
$ws.Range("J69").Value = 3746.75
$ws.Range("K69").Value = 7500
$ws.Range("L69").Value = 11240.25
$ws.Range("M69").Value = -6626
$ws.Range("N69").Value = -12988.25
$ws.Range("H72").Value = 3497.4
$ws.Range("I72").Value = 2500
$ws.Range("J72").Value = 3746.75
$ws.Range("K72").Value = 22500
$ws.Range("L72").Value = 33720.75
$ws.Range("M72").Value = -18132
$ws.Range("N72").Value = -42456.75
$ws.Range("H96").Value = 1172.25
$ws.Range("I96").Value = 1115.6
$ws.Range("J96").Value = 1266.6666
$ws.Range("K96").Value = 3346.8
$ws.Range("L96").Value = 3799.9998
$ws.Range("M96").Value = -1973.8
$ws.Range("N96").Value = -6545.9998
$ws.Range("H98").Value = 3976347.2
$ws.Range("I98").Value = 9596.546
$ws.Range("K98").Value = 9596.546
$ws.Range("M98").Value = -8098.546
$ws.Range("H100").Value = 2326.5625
$ws.Range("I100").Value = 1152.6364
$ws.Range("J100").Value = 4909.2
$ws.Range("K100").Value = 1152.6364
$ws.Range("L100").Value = 4909.2
$ws.Range("M100").Value = -611.6364000000001
$ws.Range("N100").Value = -5991.2
$ws.Range("H122").Value = 3976347.2
$ws.Range("I122").Value = 9596.546
$ws.Range("K122").Value = 28789.638
$ws.Range("M122").Value = -26339.638
$ws.Range("H129").Value = 1066.4166
$ws.Range("I129").Value = 497
$ws.Range("J129").Value = 1180.3
$ws.Range("K129").Value = 1491
$ws.Range("L129").Value = 3540.9
$ws.Range("M129").Value = 3509
$ws.Range("N129").Value = -13540.9
$ws.Range("H137").Value = 1033.7576
$ws.Range("I137").Value = 866.6070999999999
$ws.Range("J137").Value = 1969.8
$ws.Range("K137").Value = 2599.8213
$ws.Range("L137").Value = 5909.4
$ws.Range("M137").Value = -49.82129999999961
$ws.Range("N137").Value = -11009.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8200.75
$ws.Range("J37").Value = 8200.75
$ws.Range("L37").Value = 8200.75
$ws.Range("N37").Value = -8746.75
$ws.Range("H102").Value = 5425.75
$ws.Range("I102").Value = 4513.625
$ws.Range("J102").Value = 7250
$ws.Range("K102").Value = 4513.625
$ws.Range("L102").Value = 7250
$ws.Range("M102").Value = -2891.625
$ws.Range("N102").Value = -10494
$ws.Range("H132").Value = 1641.4445
$ws.Range("I132").Value = 1226.326
$ws.Range("J132").Value = 2764.7058
$ws.Range("K132").Value = 3678.978
$ws.Range("L132").Value = 8294.117400000001
$ws.Range("M132").Value = -1148.978
$ws.Range("N132").Value = -13354.1174

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2127.4666
$ws.Range("I99").Value = 1619.8889
$ws.Range("J99").Value = 2888.8333
$ws.Range("K99").Value = 1619.8889
$ws.Range("L99").Value = 2888.8333
$ws.Range("M99").Value = -121.8888999999999
$ws.Range("N99").Value = -5884.8333
$ws.Range("H105").Value = 1619.2858
$ws.Range("I105").Value = 1561.3043
$ws.Range("J105").Value = 1689.4736
$ws.Range("K105").Value = 1561.3043
$ws.Range("L105").Value = 1689.4736
$ws.Range("M105").Value = 185.6957
$ws.Range("N105").Value = -5183.4736
$ws.Range("H134").Value = 5941.5864
$ws.Range("I134").Value = 1052.2963
$ws.Range("J134").Value = 71947
$ws.Range("K134").Value = 3156.8889
$ws.Range("L134").Value = 215841
$ws.Range("M134").Value = -621.8888999999999
$ws.Range("N134").Value = -220911

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5696.7446
$ws.Range("I31").Value = 4492.2646
$ws.Range("J31").Value = 8846.923000000001
$ws.Range("K31").Value = 4492.2646
$ws.Range("L31").Value = 8846.923000000001
$ws.Range("M31").Value = -4197.2646
$ws.Range("N31").Value = -9436.923000000001
$ws.Range("H34").Value = 5696.7446
$ws.Range("I34").Value = 4492.2646
$ws.Range("J34").Value = 8846.923000000001
$ws.Range("K34").Value = 4492.2646
$ws.Range("L34").Value = 8846.923000000001
$ws.Range("M34").Value = -4290.2646
$ws.Range("N34").Value = -9250.923000000001
$ws.Range("H80").Value = 24000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 24000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H105").Value = 2280.7273
$ws.Range("I105").Value = 1681.6666
$ws.Range("K105").Value = 1681.6666
$ws.Range("M105").Value = 65.33339999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 316.25
$ws.Range("I2").Value = 56.666668
$ws.Range("J2").Value = 376.15384
$ws.Range("K2").Value = 340.000008
$ws.Range("L2").Value = 2256.92304
$ws.Range("M2").Value = -227.000008
$ws.Range("N2").Value = -2482.92304
$ws.Range("H131").Value = 5320053
$ws.Range("J131").Value = 6024946
$ws.Range("L131").Value = 18074838
$ws.Range("N131").Value = -18084918

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 2999.6667
$ws.Range("J33").Value = 2999.6667
$ws.Range("L33").Value = 2999.6667
$ws.Range("N33").Value = -3503.6667
$ws.Range("H122").Value = 2666.4285
$ws.Range("I122").Value = 2631.3157
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 7893.9471
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5443.9471
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 2900
$ws.Range("J5").Value = 2900
$ws.Range("L5").Value = 2900
$ws.Range("N5").Value = -3126
$ws.Range("H100").Value = 3128.9443
$ws.Range("I100").Value = 2403.7144
$ws.Range("J100").Value = 3590.4546
$ws.Range("K100").Value = 2403.7144
$ws.Range("L100").Value = 3590.4546
$ws.Range("M100").Value = -1862.7144
$ws.Range("N100").Value = -4672.4546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 21280222
$ws.Range("I132").Value = 27028340
$ws.Range("J132").Value = 12179.7
$ws.Range("K132").Value = 81085020
$ws.Range("L132").Value = 36539.10000000001
$ws.Range("M132").Value = -81082490
$ws.Range("N132").Value = -41599.10000000001
